$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date for every data row
# (rows 2 through 458). Update the date from 2023-10-04 (45203) to
# 2023-10-05 (45204) everywhere it currently appears.
for ($row = 2; $row -le 458; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
